$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "R_ij_t[1,3,2]"
$ws.Range("A3").Value = "R_ij_t[2,4,1]"
$ws.Range("A4").Value = "R_ij_t[3,5,2]"
$ws.Range("A5").Value = "R_ij_t[4,6,1]"
$ws.Range("A6").Value = "R_ij_t[5,1,2]"
$ws.Range("A7").Value = "R_ij_t[6,2,1]"
$ws.Range("A12").Value = "X_ij_v[3,14,5]"
$ws.Range("A17").Value = "X_ij_v[4,9,6]"
$ws.Range("A22").Value = "X_ij_v[5,18,4]"
$ws.Range("A27").Value = "X_ij_v[6,15,3]"
$ws.Range("A28").Value = "X_ij_v[7,6,3]"
$ws.Range("A29").Value = "X_ij_v[8,10,5]"
$ws.Range("A30").Value = "X_ij_v[9,13,6]"
$ws.Range("A31").Value = "X_ij_v[10,3,5]"
$ws.Range("A32").Value = "X_ij_v[11,4,6]"
$ws.Range("A33").Value = "X_ij_v[12,5,4]"
$ws.Range("A34").Value = "X_ij_v[13,11,6]"
$ws.Range("A35").Value = "X_ij_v[14,8,5]"
$ws.Range("A36").Value = "X_ij_v[15,17,3]"
$ws.Range("A37").Value = "X_ij_v[16,12,4]"
$ws.Range("A38").Value = "X_ij_v[17,7,3]"
$ws.Range("A39").Value = "X_ij_v[18,16,4]"
$ws.Range("A46").Value = "Y_c_t[7,1]"
$ws.Range("A47").Value = "Y_c_t[8,2]"
$ws.Range("A48").Value = "Y_c_t[9,1]"
$ws.Range("A49").Value = "Y_c_t[10,2]"
$ws.Range("A50").Value = "Y_c_t[11,1]"
$ws.Range("A51").Value = "Y_c_t[12,2]"
$ws.Range("A52").Value = "Y_c_t[13,1]"
$ws.Range("A53").Value = "Y_c_t[14,2]"
$ws.Range("A54").Value = "Y_c_t[15,1]"
$ws.Range("A55").Value = "Y_c_t[16,2]"
$ws.Range("A56").Value = "Y_c_t[17,1]"
$ws.Range("A57").Value = "Y_c_t[18,2]"
$ws.Range("A59").Value = "Y_c_v[8,5]"
$ws.Range("A60").Value = "Y_c_v[9,6]"
$ws.Range("A61").Value = "Y_c_v[10,5]"
$ws.Range("A62").Value = "Y_c_v[11,6]"
$ws.Range("A63").Value = "Y_c_v[12,4]"
$ws.Range("A64").Value = "Y_c_v[13,6]"
$ws.Range("A65").Value = "Y_c_v[14,5]"
$ws.Range("A67").Value = "Y_c_v[16,4]"
$ws.Range("A69").Value = "Y_c_v[18,4]"
$ws.Range("A83").Value = "B_v_s[4,5]"
$ws.Range("A84").Value = "B_v_s[5,3]"
$ws.Range("A85").Value = "B_v_s[6,4]"
$ws.Range("A86").Value = "H_ij_v[3,14,5]"
$ws.Range("A87").Value = "H_ij_v[4,9,6]"
$ws.Range("A88").Value = "H_ij_v[5,18,4]"
$ws.Range("A89").Value = "H_ij_v[6,15,3]"
$ws.Range("A90").Value = "H_ij_v[8,10,5]"
$ws.Range("B90").Value = 10
$ws.Range("A91").Value = "H_ij_v[9,13,6]"
$ws.Range("B91").Value = 20
$ws.Range("A92").Value = "H_ij_v[13,11,6]"
$ws.Range("B92").Value = 10
$ws.Range("A93").Value = "H_ij_v[14,8,5]"
$ws.Range("A94").Value = "H_ij_v[15,17,3]"
$ws.Range("A95").Value = "H_ij_v[16,12,4]"
$ws.Range("A96").Value = "H_ij_v[17,7,3]"
$ws.Range("A97").Value = "H_ij_v[18,16,4]"
$ws.Range("B97").Value = 20
$ws.Range("A98").Value = "G_i_t[5,2]"
$ws.Range("B98").Value = 1
$ws.Range("A99").Value = "G_i_t[6,1]"
$ws.Range("B99").Value = 1
$ws.Range("A100").Value = "G_i_v[7,3]"
$ws.Range("B100").Value = 2
$ws.Range("A101").Value = "G_i_v[8,5]"
$ws.Range("B101").Value = 1
$ws.Range("A102").Value = "G_i_v[10,5]"
$ws.Range("B102").Value = 2
$ws.Range("A103").Value = "G_i_v[11,6]"
$ws.Range("B103").Value = 2
$ws.Range("A104").Value = "G_i_v[12,4]"
$ws.Range("B104").Value = 2
$ws.Range("A105").Value = "G_i_v[13,6]"
$ws.Range("B105").Value = 1
$ws.Range("A106").Value = "G_i_v[16,4]"
$ws.Range("B106").Value = 1
$ws.Range("A107").Value = "G_i_v[17,3]"
$ws.Range("B107").Value = 1

# Delete rows 108-118 (no longer present in the data set)
$ws.Range("A108:B118").EntireRow.Delete()

$ws.Range("A1").Select()
